# Add ANC HIV prevalence data (commit: "add ANC HIV prevalence")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New explanatory note on the last existing "HIV Prevalence in Women 2008" row (row 43),
# in the Questions column (K).
$ws.Range("K43").Value = "Prevalence estimates were reported in ANC surveillance report 2011. The number of HIV positive women and the number of women attending ANC each year were not available, except for 2011. The number positive each year from 1990-2003 were calculated using the number of women attending ANC in 2011.  "

# New rows 44-57: "HIV prevalence in women - ANC" / "ANC Surveillance Report 2011" / "All"
# Columns: A=Criteria, B=Source, C=Group, D=Year, E=Pos (=Mean*N), F=N, G=Mean, J=Usage Status
$criteria = "HIV prevalence in women - ANC"
$source   = "ANC Surveillance Report 2011"
$group    = "All"

$data = @(
    @{Row=44; Year=1990; N=410; Mean=0.18},
    @{Row=45; Year=1991; N=410; Mean=0.18},
    @{Row=46; Year=1992; N=410; Mean=0.19},
    @{Row=47; Year=1993; N=410; Mean=0.19},
    @{Row=48; Year=1994; N=410; Mean=0.29},
    @{Row=49; Year=1995; N=410; Mean=0.24},
    @{Row=50; Year=1996; N=410; Mean=0.26},
    @{Row=51; Year=1997; N=410; Mean=0.32},
    @{Row=52; Year=1998; N=410; Mean=0.27},
    @{Row=53; Year=1999; N=410; Mean=0.25},
    @{Row=54; Year=2000; N=410; Mean=0.33},
    @{Row=55; Year=2001; N=410; Mean=0.29},
    @{Row=56; Year=2002; N=410; Mean=0.26},
    @{Row=57; Year=2003; N=410; Mean=0.26}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Range("A$r").Value = $criteria
    $ws.Range("B$r").Value = $source
    $ws.Range("C$r").Value = $group
    $ws.Range("D$r").Value = $rec.Year
    $ws.Range("F$r").Value = $rec.N
    $ws.Range("G$r").Value = $rec.Mean
    $ws.Range("E$r").Formula = "=Table1[[#This Row],[Mean]]*Table1[[#This Row],[N]]"
    $ws.Range("J$r").Value = "Y"
}

# Restore the view/selection state recorded for this edit.
$ws.Range("D16").Select()
